$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text: collapse the old blank-line-wrapped "or" separators
# in the Output/ResponseHeader columns (F, G) into a single-line "-OR-". The
# underlying route/method/input descriptions are unchanged; only this
# separator wording was revised.

$ws.Range("F2").Value = @"
Created 201
-OR-
errorCode
"@

$ws.Range("G2").Value = @"
{
user:{
id: <number>,
username: <string>,
}
}
-OR-
{
errorCode: <number>,
errorMessage: <string>
}
"@

$ws.Range("F3").Value = @"
OK 200
-OR-
errorCode
"@

$ws.Range("G3").Value = @"
{
authToken: <string>
}
-OR-
{
errorCode: <number>,
errorMessage: <string>
}
"@

$ws.Range("F4").Value = @"
OK 200
-OR-
errorCode
"@

$ws.Range("G4").Value = @"
none
-OR-
{
errorCode: <number>,
errorMessage: <string>
}
"@

$ws.Range("F5").Value = @"
OK 200
-OR-
errorCode
"@

$ws.Range("F6").Value = @"
OK 200 (No error, but old value was higher),
CREATED 201 (New highscore stored)
-OR-
errorCode
"@

$ws.Range("F7").Value = @"
OK 200
-OR-
errorCode
"@

$ws.Range("G7").Value = @"
{
'<number (levelIndex)>' : <number> (highscore)
}
-OR-
{
errorCode: <number>,
errorMessage: <string>
}
"@

$ws.Range("F8").Value = @"
OK 200,
-OR-
errorCode
"@

# --- Row heights shrink now that the blank lines around "or" are gone
# (content is auto-wrapped; the rows that lost text height need to be
# resized to match the new line counts).
$ws.Rows.Item(2).RowHeight = 165
$ws.Rows.Item(3).RowHeight = 120
$ws.Rows.Item(4).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(7).RowHeight = 150

# --- View state: zoomed to 75% with H7 selected (previously topLeftCell
# A7 / D12 selected).
$win = $excel.ActiveWindow
$win.Zoom = 75
$ws.Range("H7").Select()
